$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Actual Result" (H) / "Test Status" (I) need to be reset
# to the "not executed" state. This also normalizes the I-column
# formatting (PASSED=green / FAILED=red) back to the plain/no-fill
# style already used by column H, by copying H's format onto I.
$rows = 6,7,8,9,10,11,12,13,14,15,16

foreach ($r in $rows) {
    $hCell = $ws.Range("H$r")
    $iCell = $ws.Range("I$r")

    $hCell.Value = "Test not executed"
    $iCell.Value = "Not Run"

    # Re-use H's existing cell style (no fill) for I instead of the
    # PASSED/FAILED colored style it had before.
    $hCell.Copy()
    $iCell.PasteSpecial(-4122)
}
